# sem3_timetable.xlsx update:
#  - rebalance the elective-basket scheduling in Section_A / Section_B
#    (common slots for CS261/CS263/CS264/MA261, incl. a CS264 Tutorial slot)
#  - add a new "Course_Summary" sheet listing course codes/names/credits/instructors

$wb = $excel.ActiveWorkbook
$wsA = $wb.Worksheets.Item("Section_A")
$wsB = $wb.Worksheets.Item("Section_B")

# ---- Section_A (sheet1) timetable changes ----
$wsA.Range("D2").Value = "CS261"
$wsA.Range("E2").Value = "CS264"
$wsA.Range("F2").Value = "Free"

$wsA.Range("C3").Value = "Free"
$wsA.Range("D3").Value = "Free"
$wsA.Range("F3").Value = "CS261"

$wsA.Range("B5").Value = "Free"
$wsA.Range("C5").Value = "CS261"
$wsA.Range("E5").Value = "CS263"
$wsA.Range("F5").Value = "Free"

$wsA.Range("B6").Value = "MA261"
$wsA.Range("D6").Value = "CS264"
$wsA.Range("E6").Value = "CS264 (Tutorial)"
$wsA.Range("F6").Value = "CS264"

$wsA.Range("C7").Value = "MA261"
$wsA.Range("D7").Value = "Free"
$wsA.Range("F7").Value = "Free"

# ---- Section_B (sheet2) timetable changes ----
$wsB.Range("B2").Value = "Free"
$wsB.Range("C2").Value = "CS264"
$wsB.Range("E2").Value = "CS263"
$wsB.Range("F2").Value = "CS264 (Tutorial)"

$wsB.Range("B3").Value = "Free"
$wsB.Range("D3").Value = "CS261"
$wsB.Range("E3").Value = "MA261"
$wsB.Range("F3").Value = "Free"

$wsB.Range("B5").Value = "CS264"
$wsB.Range("C5").Value = "CS263"
$wsB.Range("D5").Value = "MA261"

$wsB.Range("C6").Value = "Free"
$wsB.Range("D6").Value = "CS264"
$wsB.Range("E6").Value = "Free"
$wsB.Range("F6").Value = "Free"

$wsB.Range("C7").Value = "CS261"

# ---- New Course_Summary sheet (added after Section_B) ----
$wsSummary = $wb.Worksheets.Add($null, $wsB)
$wsSummary.Name = "Course_Summary"

$wsSummary.Range("A1").Value = "Course Code"
$wsSummary.Range("B1").Value = "Course Name"
$wsSummary.Range("C1").Value = "Course Type"
$wsSummary.Range("D1").Value = "LTPSC"
$wsSummary.Range("E1").Value = "Credits"
$wsSummary.Range("F1").Value = "Instructor"

$wsA.Range("B1").Copy()
$wsSummary.Range("A1:F1").PasteSpecial(-4122)

$wsSummary.Range("A2").Value = "MA261"
$wsSummary.Range("B2").Value = "Differential Equations"
$wsSummary.Range("C2").Value = "Core"
$wsSummary.Range("D2").Value = "2-0-0-0-2"
$wsSummary.Range("E2").Value = 2
$wsSummary.Range("F2").Value = "Dr. Anand Barangi"

$wsSummary.Range("A3").Value = "CS261"
$wsSummary.Range("B3").Value = "Operating System"
$wsSummary.Range("C3").Value = "Core"
$wsSummary.Range("D3").Value = "3-0-0-4-2"
$wsSummary.Range("E3").Value = 5
$wsSummary.Range("F3").Value = "Dr. Somes"

$wsSummary.Range("A4").Value = "CS263"
$wsSummary.Range("B4").Value = "Design & Analysis of Algorithms"
$wsSummary.Range("C4").Value = "Core"
$wsSummary.Range("D4").Value = "3-0-2-0-4"
$wsSummary.Range("E4").Value = 7
$wsSummary.Range("F4").Value = "Dr. Prabhu Prasad"

$wsSummary.Range("A5").Value = "CS264"
$wsSummary.Range("B5").Value = "Computer Networks"
$wsSummary.Range("C5").Value = "Core"
$wsSummary.Range("D5").Value = "3-1-0-0-4"
$wsSummary.Range("E5").Value = 8
$wsSummary.Range("F5").Value = "Dr. Prabhu Prasad"
